$d = $word.ActiveDocument

# The 5x5 table has equations only on rows 1, 5, 10, 15, 20 (the other
# rows are blank spacer rows). We walk the table in row-major order and
# replace the text of every non-empty equation cell with the new value,
# preserving document order so duplicate equations (e.g. "87x52=4524")
# are each mapped to their correct, distinct replacement.

$newValues = @(
    "41x52=2132",
    "39x47=1833",
    "45x62=2790",
    "71x32=2272",
    "72x16=1152",
    "65x60=3900",
    "32x54=1728",
    "49x52=2548",
    "51x16=816",
    "45x95=4275",
    "41x79=3239",
    "50x77=3850",
    "18x42=756",
    "41x79=3239",
    "91x48=4368",
    "87x95=8265",
    "35x58=2030",
    "14x55=770",
    "19x54=1026",
    "83x59=4897",
    "78x96=7488",
    "27x69=1863",
    "26x28=728",
    "49x33=1617",
    "38x33=1254"
)
# multiplication sign used in the document
$times = [char]0x00D7
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $newValues[$i] = $newValues[$i].Replace("x", $times)
}

$tbl = $d.Tables(1)
$idx = 0
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
        $cell = $tbl.Cell($r, $c)
        $rng = $cell.Range
        # trim the trailing cell-mark characters to get the visible text
        $txt = $rng.Text
        $txt = $txt.TrimEnd([char]7).TrimEnd([char]13)
        if ($txt.Length -gt 0) {
            $rng.Text = $newValues[$idx]
            $idx = $idx + 1
        }
    }
}

Write-Output "replaced $idx cells"
